$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.430.35'
$ws.Range("E2").Value = '  +1.66%  '

$ws.Range("D3").Value = '1.862.95'
$ws.Range("E3").Value = '  +0.85%  '

$ws.Range("D4").Value = '1.011'
$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").Value = '310.71'
$ws.Range("E5").Value = '  +0.23%  '

$ws.Range("D6").Value = '1.010'

$ws.Range("D7").Value = '0.4778'
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '0.3762'
$ws.Range("E8").Value = '  +2.20%  '

$ws.Range("D9").Value = '0.07328'
$ws.Range("E9").Value = '  +1.34%  '

$ws.Range("D10").Value = '0.9343'
$ws.Range("E10").Value = '  +0.47%  '

$ws.Range("E11").Value = '  +4.56%  '

$ws.Range("E12").Value = '  +1.20%  '

$ws.Range("D13").Value = '1.905.20'
$ws.Range("E13").Value = '  +4.33%  '

$ws.Range("D14").Value = '5.428'
$ws.Range("E14").Value = '  +1.73%  '

$ws.Range("D15").Value = '6.554'
$ws.Range("E15").Value = '  +1.90%  '

$ws.Range("D16").Value = '90.37'
$ws.Range("E16").Value = '  +1.80%  '

$ws.Range("E17").Value = '  -0.28%  '

$ws.Range("D18").Value = '0.000008879'
$ws.Range("E18").Value = '  +2.79%  '

$ws.Range("E19").Value = '  -0.18%  '

$ws.Range("D20").Value = '27.498.93'
$ws.Range("E20").Value = '  +1.76%  '

$ws.Range("D21").Value = '14.70'
$ws.Range("E21").Value = '  +1.48%  '

$ws.Range("D22").Value = '5.112'
$ws.Range("E22").Value = '  +1.08%  '

$ws.Range("E23").Value = '  +0.48%  '

$ws.Range("D24").Value = '1.935'
$ws.Range("E24").Value = '  +0.36%  '

$ws.Range("D25").Value = '155.48'
$ws.Range("E25").Value = '  +1.79%  '

$ws.Range("D26").Value = '18.47'
$ws.Range("E26").Value = '  +1.39%  '

$ws.Range("D27").Value = '2.019'
$ws.Range("E27").Value = '  +1.00%  '

$ws.Range("D28").Value = '115.37'
$ws.Range("E28").Value = '  +0.89%  '

$ws.Range("E29").Value = '  -0.85%  '

$ws.Range("D30").Value = '0.08893'
$ws.Range("E30").Value = '  -0.09%  '

$ws.Range("E32").Value = '  +3.27%  '

$ws.Range("D33").Value = '0.7544'
$ws.Range("E33").Value = '  +1.51%  '

$ws.Range("D34").Value = '4.599'
$ws.Range("E34").Value = '  +2.11%  '

$ws.Range("D35").Value = '2.746'
$ws.Range("E35").Value = '  +0.32%  '

$ws.Range("D36").Value = '0.02052'
$ws.Range("E36").Value = '  +4.83%  '

$ws.Range("E37").Value = '  -0.18%  '

$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").Value = '0.5492'
$ws.Range("E38").Value = '  +5.22%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.05268'
$ws.Range("E39").Value = '  +0.10%  '

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.986'
$ws.Range("E40").Value = '  +0.30%  '

$ws.Range("D41").Value = '7.063'
$ws.Range("E41").Value = '  +1.02%  '

$ws.Range("D42").Value = '8.699'
$ws.Range("E42").Value = '  +6.04%  '

$ws.Range("E43").Value = '  +0.72%  '

$ws.Range("D44").Value = '10.69'
$ws.Range("E44").Value = '  +0.40%  '

$ws.Range("D45").Value = '0.4852'
$ws.Range("E45").Value = '  +2.17%  '

$ws.Range("E46").Value = '  -0.12%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '103.09'
$ws.Range("E47").Value = '  +1.30%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '1.656'
$ws.Range("E48").Value = '  +2.86%  '

$ws.Range("D49").Value = '67.32'
$ws.Range("E49").Value = '  +2.68%  '

$ws.Range("E50").Value = '  +0.30%  '

$ws.Range("D51").Value = '0.9183'
$ws.Range("E51").Value = '  +3.53%  '
